$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5037.4326
$ws.Range("I40").Value = 2937
$ws.Range("J40").Value = 5365.625
$ws.Range("K40").Value = 2937
$ws.Range("L40").Value = 5365.625
$ws.Range("M40").Value = -2762
$ws.Range("N40").Value = -5715.625

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2999.9658
$ws.Range("I51").Value = 2999.9727
$ws.Range("J51").Value = 2999.9333
$ws.Range("K51").Value = 2999.9727
$ws.Range("L51").Value = 2999.9333
$ws.Range("M51").Value = -2515.9727
$ws.Range("N51").Value = -3967.9333

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2279.9412
$ws.Range("J112").Value = 2366
$ws.Range("L112").Value = 7098
$ws.Range("N112").Value = -9314

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 47378.285
$ws.Range("I116").Value = 45250
$ws.Range("J116").Value = 47733
$ws.Range("K116").Value = 45250
$ws.Range("L116").Value = 47733
$ws.Range("M116").Value = -41808
$ws.Range("N116").Value = -54617

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5755.3125
$ws.Range("J138").Value = 5206.2856
$ws.Range("L138").Value = 15618.8568
$ws.Range("N138").Value = -25898.8568

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1323
$ws.Range("J2").Value = 1196
$ws.Range("L2").Value = 1196
$ws.Range("N2").Value = -1422

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2210.6
$ws.Range("J45").Value = 3183.3333
$ws.Range("L45").Value = 3183.3333
$ws.Range("N45").Value = -3937.3333

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2263.7917
$ws.Range("I61").Value = 2279.05
$ws.Range("J61").Value = 2187.5
$ws.Range("K61").Value = 2279.05
$ws.Range("L61").Value = 2187.5
$ws.Range("M61").Value = -2067.05
$ws.Range("N61").Value = -2611.5

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1544.3636
$ws.Range("I102").Value = 1544.3636
$ws.Range("K102").Value = 1544.3636
$ws.Range("M102").Value = 77.63640000000009

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1454.5
$ws.Range("I110").Value = 1411.8334
$ws.Range("J110").Value = 1497.1666
$ws.Range("K110").Value = 1411.8334
$ws.Range("L110").Value = 1497.1666
$ws.Range("M110").Value = 633.1666
$ws.Range("N110").Value = -5587.1666

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1323
$ws.Range("J116").Value = 1196
$ws.Range("L116").Value = 1196
$ws.Range("N116").Value = -5784

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5037.271
$ws.Range("I132").Value = 2947.575
$ws.Range("K132").Value = 8842.724999999999
$ws.Range("M132").Value = -6312.724999999999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2263.7917
$ws.Range("I136").Value = 2279.05
$ws.Range("J136").Value = 2187.5
$ws.Range("K136").Value = 6837.150000000001
$ws.Range("L136").Value = 6562.5
$ws.Range("M136").Value = -4287.150000000001
$ws.Range("N136").Value = -11662.5

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1323
$ws.Range("J3").Value = 1196
$ws.Range("L3").Value = 1196
$ws.Range("N3").Value = -1424

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2873
$ws.Range("I86").Value = 2800
$ws.Range("J86").Value = 3128.5
$ws.Range("K86").Value = 2800
$ws.Range("L86").Value = 3128.5
$ws.Range("M86").Value = -1677
$ws.Range("N86").Value = -5374.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2873
$ws.Range("I89").Value = 2800
$ws.Range("J89").Value = 3128.5
$ws.Range("K89").Value = 14000
$ws.Range("L89").Value = 15642.5
$ws.Range("M89").Value = -8384
$ws.Range("N89").Value = -26874.5

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1165.4546
$ws.Range("I107").Value = 1165.4546
$ws.Range("K107").Value = 1165.4546
$ws.Range("M107").Value = 754.5454

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2045.7446
$ws.Range("I134").Value = 1790.7646
$ws.Range("K134").Value = 5372.293799999999
$ws.Range("M134").Value = -2837.293799999999

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125510.555
$ws.Range("I16").Value = 25200
$ws.Range("J16").Value = 250898.75
$ws.Range("K16").Value = 25200
$ws.Range("L16").Value = 250898.75
$ws.Range("M16").Value = -24913
$ws.Range("N16").Value = -251472.75

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1426.4348
$ws.Range("I22").Value = 738.75
$ws.Range("J22").Value = 2998.2856
$ws.Range("K22").Value = 738.75
$ws.Range("L22").Value = 2998.2856
$ws.Range("M22").Value = -388.75
$ws.Range("N22").Value = -3698.2856

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2048.6553
$ws.Range("I31").Value = 1836.1072
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 1836.1072
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -1541.1072
$ws.Range("N31").Value = -8590

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2048.6553
$ws.Range("I34").Value = 1836.1072
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 1836.1072
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -1634.1072
$ws.Range("N34").Value = -8404

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4940.5
$ws.Range("I58").Value = 3383
$ws.Range("K58").Value = 3383
$ws.Range("M58").Value = -3180

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2864.3333
$ws.Range("I99").Value = 2705.625
$ws.Range("K99").Value = 2705.625
$ws.Range("M99").Value = -1207.625

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 125510.555
$ws.Range("I113").Value = 25200
$ws.Range("J113").Value = 250898.75
$ws.Range("K113").Value = 25200
$ws.Range("L113").Value = 250898.75
$ws.Range("M113").Value = -23030
$ws.Range("N113").Value = -255238.75

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2864.3333
$ws.Range("I126").Value = 2705.625
$ws.Range("K126").Value = 8116.875
$ws.Range("M126").Value = -5646.875

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2505.7058
$ws.Range("I132").Value = 2391.724
$ws.Range("K132").Value = 7175.172
$ws.Range("M132").Value = -4645.172

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2161.5264
$ws.Range("I134").Value = 1843.7
$ws.Range("K134").Value = 5531.1
$ws.Range("M134").Value = -2996.1

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4940.5
$ws.Range("I136").Value = 3383
$ws.Range("K136").Value = 10149
$ws.Range("M136").Value = -7599

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12003306
$ws.Range("I4").Value = 20000844
$ws.Range("K4").Value = 60002532
$ws.Range("M4").Value = -60002420

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 166.38461
$ws.Range("I23").Value = 81.75
$ws.Range("J23").Value = 204
$ws.Range("K23").Value = 245.25
$ws.Range("L23").Value = 612
$ws.Range("M23").Value = -10.25
$ws.Range("N23").Value = -1082

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 95000
$ws.Range("J37").Value = 95000
$ws.Range("L37").Value = 285000
$ws.Range("N37").Value = -285224

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 428.66666
$ws.Range("I97").Value = 297.16666
$ws.Range("J97").Value = 691.6667
$ws.Range("K97").Value = 891.4999799999999
$ws.Range("L97").Value = 2075.0001
$ws.Range("M97").Value = -395.4999799999999
$ws.Range("N97").Value = -3067.0001

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1490.6428
$ws.Range("I132").Value = 1874.3334
$ws.Range("K132").Value = 16869.0006
$ws.Range("M132").Value = -14339.0006

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 6465.1113
$ws.Range("I141").Value = 5023.25
$ws.Range("K141").Value = 15069.75
$ws.Range("M141").Value = -9889.75

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2137.037
$ws.Range("I122").Value = 1466.8334
$ws.Range("J122").Value = 3477.4443
$ws.Range("K122").Value = 4400.5002
$ws.Range("L122").Value = 10432.3329
$ws.Range("M122").Value = -1950.5002
$ws.Range("N122").Value = -15332.3329

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2775.16
$ws.Range("I22").Value = 2291.8462
$ws.Range("K22").Value = 2291.8462
$ws.Range("M22").Value = -1996.8462

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2775.16
$ws.Range("I27").Value = 2291.8462
$ws.Range("K27").Value = 2291.8462
$ws.Range("M27").Value = -2184.8462

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1995.125
$ws.Range("I100").Value = 2123
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 2123
$ws.Range("L100").Value = 1100
$ws.Range("M100").Value = -1582
$ws.Range("N100").Value = -2182

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2426.2812
$ws.Range("I132").Value = 1538.6316
$ws.Range("J132").Value = 3723.6155
$ws.Range("K132").Value = 2426.2812
$ws.Range("L132").Value = 11170.8465
$ws.Range("M132").Value = -2085.8948
$ws.Range("N132").Value = -16230.8465

# WVR row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 787500
$ws.Range("I21").Value = 1500000
$ws.Range("J21").Value = 75000
$ws.Range("K21").Value = 1500000
$ws.Range("L21").Value = 75000
$ws.Range("M21").Value = -1499765
$ws.Range("N21").Value = -75470

# WVR row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 787500
$ws.Range("I35").Value = 1500000
$ws.Range("J35").Value = 75000
$ws.Range("K35").Value = 1500000
$ws.Range("L35").Value = 75000
$ws.Range("M35").Value = -1499710
$ws.Range("N35").Value = -75580

# WVR row 49
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 369300
$ws.Range("J49").Value = 53950
$ws.Range("L49").Value = 53950
$ws.Range("N49").Value = -54410

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1750
$ws.Range("I96").Value = 1750
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1750
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -377
$ws.Range("N96").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2366.72
$ws.Range("J132").Value = 3581.25
$ws.Range("L132").Value = 10743.75
$ws.Range("N132").Value = -15803.75
